# Updated cryptos list on Mon Sep 30 07:43:44 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text rather than
# reinterpreting number-like strings (e.g. "68.30") as numeric values.

$ws.Range("D2").Value = "'64.567.20"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "'2.637.39"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'582.79"
$ws.Range("E5").Value = "  -2.20%  "

$ws.Range("D6").Value = "'157.17"
$ws.Range("E6").Value = "  +1.04%  "

$ws.Range("E7").Value = "  +2.53%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -2.82%  "

$ws.Range("D10").Value = "'5.83"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "'28.72"
$ws.Range("E13").Value = "  +0.17%  "

$ws.Range("D14").Value = "'0.0000187"
$ws.Range("E14").Value = "  -4.29%  "

$ws.Range("D15").Value = "'3.115.73"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "'64.383.24"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").Value = "'2.634.20"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "'12.26"
$ws.Range("E18").Value = "  -2.61%  "

$ws.Range("D19").Value = "'4.70"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("D20").Value = "'7.48"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").Value = "'348.02"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'68.30"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").Value = "'1.77"
$ws.Range("E24").Value = "  +7.78%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "'9.45"
$ws.Range("E26").Value = "  -1.38%  "

$ws.Range("D27").Value = "'594.71"
$ws.Range("E27").Value = "  +10.42%  "

$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "'8.02"
$ws.Range("E29").Value = "  +1.82%  "

$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("D34").Value = "'6.68"
$ws.Range("E34").Value = "  +4.97%  "

$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").Value = "'20.09"
$ws.Range("E37").Value = "  -0.79%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("D40").Value = "'153.65"
$ws.Range("E40").Value = "  -0.77%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  +5.16%  "

$ws.Range("D43").Value = "'158.80"
$ws.Range("E43").Value = "  -0.92%  "

$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("D45").Value = "'23.47"
$ws.Range("E45").Value = "  +4.41%  "

$ws.Range("D46").Value = "'0.0606"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").Value = "'0.637"
$ws.Range("E47").Value = "  +0.48%  "

$ws.Range("D48").Value = "'0.0257"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("D50").Value = "'19.27"
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("D51").Value = "'0.0₆0238"
$ws.Range("E51").Value = "  -5.67%  "
